# Laborator 12.02.2024 - Matrici
# Adds a new week's attendance (column D) for existing students, enrolls
# six new students (rows 26-31) with this week's attendance, then
# re-sorts the roster (B3:S31) alphabetically by name - mirroring the
# original spreadsheet's "sorted" attendance table workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Mark this week's attendance (column D) for the students already on
#    the roster (rows 3-25). Students that attended get D = TRUE; the
#    rest are left blank (no D cell at all), same convention as column C.
# ---------------------------------------------------------------------
$attendedThisWeek = @(3,5,6,9,10,11,13,15,17,19,20,23,24)
foreach ($r in $attendedThisWeek) {
    $ws.Range("D$r").Value = $true
}

# ---------------------------------------------------------------------
# 2) Enroll the six new students into the previously empty rows
#    (26-31). They were not present in week 1 (column C left blank),
#    but did attend this week (column D = TRUE).
# ---------------------------------------------------------------------
$newStudents = @(
    "Aniko Vieriu",
    "Diana Lazea",
    "Cosmin Chira",
    "Vlad Chis",
    "Patrick Tocut",
    "Alexandru Lupse"
)
$r = 26
foreach ($name in $newStudents) {
    $ws.Range("B$r").Value = $name
    $ws.Range("D$r").Value = $true
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Re-sort the full roster (B3:S31) alphabetically by first name,
#    exactly like the pre-existing sortState on the sheet.
# ---------------------------------------------------------------------
$sortRange = $ws.Range("B3:S31")
$keyRange = $ws.Range("B3:B31")
$sortRange.Sort($keyRange, 1)

# ---------------------------------------------------------------------
# 4) Refresh the view: new zoom level and active selection.
# ---------------------------------------------------------------------
$sheetView = $ws.Application.ActiveWindow
$sheetView.Zoom = 160
$sheetView.ScrollRow = 1
$ws.Range("H14").Select()
